# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the zh-cn / de-de
# handoffs have moved from "In Translation" to "Ready for handoff", and
# refreshes the associated timestamps. Also widens the Status columns to
# fit the longer "Ready for handoff" text (mirrors Excel's column autosize
# after the content grows).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
# ColumnWidth value whose stored (post round-trip) width lands on the
# target ~17.22 character-units used throughout the workbook.
$newColWidth = 16.38265482584637

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-24 02:58:05"
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-24 02:57:57"
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-24 02:58:05"
$dede.Columns.Item(3).ColumnWidth = $newColWidth
